# Add a new "2022-Q3" worksheet (fund-holdings detail) right after "总计",
# before the existing "2022-Q2" sheet, and add a matching summary row on
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet
#    (keeps header/index-column styling identical) and inserting it
#    right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The copied sheet has 12 data rows (rows 2-12, plus header row 1);
# the 2022-Q3 data only has 7 funds (rows 2-8) so drop rows 9-12.
$q3Sheet.Range("A9:H12").EntireRow.Delete()

# Force the fund-code / numeric-text columns (B:G) to Text format so the
# overwritten values keep their original string representation (leading
# zeros in fund codes, "x.xx" formatted numbers) instead of being
# auto-coerced to numbers.
$q3Sheet.Range("B2:G8").NumberFormat = "@"

# Row 2: 090016 / 大成消费主题混合
$q3Sheet.Range("B2").Value = "090016"
$q3Sheet.Range("C2").Value = "大成消费主题混合"
$q3Sheet.Range("D2").Value = "3.59"
$q3Sheet.Range("E2").Value = "90.85"
$q3Sheet.Range("F2").Value = "5.86"
$q3Sheet.Range("G2").Value = "0.2104"
$q3Sheet.Range("H2").NumberFormat = "General"
$q3Sheet.Range("H2").Value = 4

# Row 3: 001195 / 工银农业产业股票
$q3Sheet.Range("B3").Value = "001195"
$q3Sheet.Range("C3").Value = "工银农业产业股票"
$q3Sheet.Range("D3").Value = "5.53"
$q3Sheet.Range("E3").Value = "80.67"
$q3Sheet.Range("F3").Value = "3.32"
$q3Sheet.Range("G3").Value = "0.1836"
$q3Sheet.Range("H3").NumberFormat = "General"
$q3Sheet.Range("H3").Value = 7

# Row 4: 002319 / 大成一带一路灵活配置混合
$q3Sheet.Range("B4").Value = "002319"
$q3Sheet.Range("C4").Value = "大成一带一路灵活配置混合"
$q3Sheet.Range("D4").Value = "1.26"
$q3Sheet.Range("E4").Value = "89.65"
$q3Sheet.Range("F4").Value = "4.26"
$q3Sheet.Range("G4").Value = "0.0537"
$q3Sheet.Range("H4").NumberFormat = "General"
$q3Sheet.Range("H4").Value = 10

# Row 5: 160323 / 华夏磐泰混合（LOF）A
$q3Sheet.Range("B5").Value = "160323"
$q3Sheet.Range("C5").Value = "华夏磐泰混合（LOF）A"
$q3Sheet.Range("D5").Value = "6.20"
$q3Sheet.Range("E5").Value = "28.65"
$q3Sheet.Range("F5").Value = "0.56"
$q3Sheet.Range("G5").Value = "0.0347"
$q3Sheet.Range("H5").NumberFormat = "General"
$q3Sheet.Range("H5").Value = 7

# Row 6: 013360 / 华夏磐泰混合（LOF）C
$q3Sheet.Range("B6").Value = "013360"
$q3Sheet.Range("C6").Value = "华夏磐泰混合（LOF）C"
$q3Sheet.Range("D6").Value = "3.70"
$q3Sheet.Range("E6").Value = "28.65"
$q3Sheet.Range("F6").Value = "0.56"
$q3Sheet.Range("G6").Value = "0.0207"
$q3Sheet.Range("H6").NumberFormat = "General"
$q3Sheet.Range("H6").Value = 7

# Row 7: 161038 / 富国新兴成长量化精选混合（LOF）A
$q3Sheet.Range("B7").Value = "161038"
$q3Sheet.Range("C7").Value = "富国新兴成长量化精选混合（LOF）A"
$q3Sheet.Range("D7").Value = "0.81"
$q3Sheet.Range("E7").Value = "91.98"
$q3Sheet.Range("F7").Value = "1.60"
$q3Sheet.Range("G7").Value = "0.0130"
$q3Sheet.Range("H7").NumberFormat = "General"
$q3Sheet.Range("H7").Value = 4

# Row 8: 014171 / 富国新兴成长量化精选混合（LOF）C (market value is a bare 0)
$q3Sheet.Range("B8").Value = "014171"
$q3Sheet.Range("C8").Value = "富国新兴成长量化精选混合（LOF）C"
$q3Sheet.Range("D8").Value = "0.00"
$q3Sheet.Range("E8").Value = "91.98"
$q3Sheet.Range("F8").Value = "1.60"
$q3Sheet.Range("G8").NumberFormat = "General"
$q3Sheet.Range("G8").Value = 0
$q3Sheet.Range("H8").NumberFormat = "General"
$q3Sheet.Range("H8").Value = 4

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q3" row to the "总计" summary sheet.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Row-insert can bleed formatting in from neighbouring rows; start the
# new row from a clean slate before (re)applying the look we want.
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.52

# Match the bold / centered / thin-border look used by the other
# index cells in column A.
$totalSheet.Range("A2").Borders.LineStyle = 1
$totalSheet.Range("A2").Font.Bold = $true
$totalSheet.Range("A2").HorizontalAlignment = -4108
$totalSheet.Range("A2").VerticalAlignment = -4160
